# Apply the "adding term 2.0.0 with four more concepts in FBOE" edit:
#  - Metadata sheet: bump Version, Date, and Contact display text.
#  - "Include from SNOMED CT" sheet: replace the first concept code and
#    drop the second concept row (116223007) entirely, shifting the
#    trailing rows up.

$wb = $excel.ActiveWorkbook

# ----- Metadata sheet -----
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-04T14:59:10+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# ----- "Include from SNOMED CT" sheet -----
$snomed = $wb.Worksheets.Item("Include from SNOMED CT")

# Replace the concept code in A2 while keeping it stored as TEXT (not a
# number) and keeping the existing cell style: build it as a formula that
# evaluates to the text, then paste the computed value back over itself
# (Copy + PasteSpecial values-only), which mimics Excel's own
# "convert formula to static value" idiom without touching NumberFormat.
$cell = $snomed.Range("A2")
$cell.Formula = '="116224001"'
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues
$snomed.Application.CutCopyMode = $false

# Remove the old second concept row (was 116223007); rows below shift up.
$snomed.Rows.Item(3).Delete()
